$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "hockey knee pads",
    "knee black leggings",
    "knee compression leggings",
    "knee compression pad",
    "knee compression running",
    "knee gel pad",
    "knee pad for men",
    "knee pad for volleyball",
    "knee pad for workout",
    "knee pad leggings",
    "knee pads athletic",
    "knee pads basketball youth",
    "knee pads black",
    "knee pads boys",
    "knee pads extra large",
    "knee pads for football",
    "knee pads for gym men",
    "knee pads for men",
    "knee pads for running men",
    "knee pads for soccer",
    "knee pads honeycomb",
    "knee pads mens",
    "knee pads skating youth",
    "knee pads snowboarding",
    "knee pads soccer",
    "knee pads sport",
    "knee pads squat",
    "knee pads thick",
    "knee pads weightlifting",
    "knee pads youth",
    "knee pants for men",
    "knee protection",
    "knee protector gym",
    "knee protectors",
    "knee replacement aids",
    "knee support gym",
    "knee support pants",
    "knee support pants men",
    "kneepads for volleyball",
    "leg compression pants",
    "leg compression tights",
    "leg protectors for men",
    "legging for basketball boys",
    "leggings for men sport gym",
    "leggings knee",
    "leggings medium",
    "leggings mens",
    "leggins training",
    "lightweight athletic pants for men",
    "lightweight pants men",
    "lightweight sports pants men",
    "mcdavid basketball knee pads 6446",
    "mcdavid basketball knee pads black",
    "men basketball tights",
    "men capri pants",
    "men compression pants",
    "men leggings pack",
    "men leggings running",
    "men leggings tall",
    "men running tights nike",
    "men tights legging",
    "mens athletic compression pants",
    "mens athletic leggings black",
    "mens basketball knee pads",
    "mens compression 3 4 pants",
    "mens compression leggings 3 4",
    "mens compression tights pants",
    "mens cycling tights",
    "mens hiking pants lightweight",
    "mens jogging tights",
    "mens lacrosse pads",
    "mens leggings",
    "mens lightweight workout pants",
    "mens running tights green",
    "mens running tights orange",
    "mens running tights yellow",
    "mens soccer clothing",
    "mens stretch thermal pants",
    "mens swim leggings",
    "mens tights basketball",
    "mens tights capri",
    "mens tights leggings",
    "mens underarmour snow pants",
    "mens workout pants",
    "multicam pants with knee pads",
    "nike basketball pads",
    "nike volleyball knee pads youth girls",
    "padded compression",
    "padded soccer pants",
    "pain in thigh joint",
    "pant knee pad inserts",
    "pantalon con rodilleras",
    "pants with knee pads",
    "polyester pants men pants",
    "protective knee pads for men",
    "rash guard men bjj",
    "reebok compression pants",
    "running compression leg",
    "running knee",
    "running pants youth"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
